$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 899.5
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568

$ws.Range("H51").Value = 93333.164
$ws.Range("I51").Value = 12000
$ws.Range("K51").Value = 12000
$ws.Range("M51").Value = -11516

$ws.Range("H70").Value = 3500
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 3500
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H86").Value = 59502
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 59502
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 59502
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -61748

$ws.Range("H89").Value = 59502
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 59502
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 297510
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -308742

$ws.Range("H100").Value = 725
$ws.Range("I100").Value = 616
$ws.Range("J100").Value = 997.5
$ws.Range("K100").Value = 616
$ws.Range("L100").Value = 997.5
$ws.Range("M100").Value = -75
$ws.Range("N100").Value = -2079.5

$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").ClearContents()

$ws.Range("H106").Value = 500010000
$ws.Range("I106").Value = 500010000
$ws.Range("K106").Value = 500010000
$ws.Range("M106").Value = -500009369

$ws.Range("H112").Value = 2463.4146
$ws.Range("J112").Value = 2463.4146
$ws.Range("L112").Value = 7390.2438
$ws.Range("N112").Value = -9606.2438

$ws.Range("H137").Value = 767
$ws.Range("I137").Value = 767
$ws.Range("K137").Value = 2301
$ws.Range("M137").Value = 249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 11748.25
$ws.Range("J41").Value = 14666
$ws.Range("L41").Value = 14666
$ws.Range("N41").Value = -15494

$ws.Range("H45").Value = 1695.5
$ws.Range("I45").Value = 1578.3611
$ws.Range("J45").Value = 2749.75
$ws.Range("K45").Value = 1578.3611
$ws.Range("L45").Value = 2749.75
$ws.Range("M45").Value = -1201.3611
$ws.Range("N45").Value = -3503.75

$ws.Range("H97").Value = 16667949
$ws.Range("I97").Value = 23810584
$ws.Range("J97").Value = 1799
$ws.Range("K97").Value = 23810584
$ws.Range("L97").Value = 1799
$ws.Range("M97").Value = -23810088
$ws.Range("N97").Value = -2791

$ws.Range("H122").Value = 6458.1665
$ws.Range("I122").Value = 8725
$ws.Range("K122").Value = 26175
$ws.Range("M122").Value = -23725

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 13344.5
$ws.Range("J15").Value = 13344.5
$ws.Range("L15").Value = 13344.5
$ws.Range("N15").Value = -13684.5

$ws.Range("H31").Value = 1638.4
$ws.Range("I31").Value = 1298
$ws.Range("K31").Value = 1298
$ws.Range("M31").Value = -1003

$ws.Range("H34").Value = 1638.4
$ws.Range("I34").Value = 1298
$ws.Range("K34").Value = 1298
$ws.Range("M34").Value = -1096

$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51498

$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -157488

$ws.Range("H105").Value = 1996
$ws.Range("I105").Value = 1996
$ws.Range("K105").Value = 1996
$ws.Range("M105").Value = -249

$ws.Range("H122").Value = 9666.333000000001
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 13499.5
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 40498.5
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -45398.5

$ws.Range("H134").Value = 3198
$ws.Range("I134").Value = 2064.3333
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 6192.999899999999
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -3657.999899999999
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 4
$ws.Range("I47").Value = 4
$ws.Range("K47").Value = 12
$ws.Range("M47").Value = 419

$ws.Range("H50").Value = 1580.1875
$ws.Range("I50").Value = 397.57144
$ws.Range("K50").Value = 1192.71432
$ws.Range("M50").Value = -711.71432

$ws.Range("H53").Value = 1580.1875
$ws.Range("I53").Value = 397.57144
$ws.Range("K53").Value = 1192.71432
$ws.Range("M53").Value = -711.71432

$ws.Range("H68").Value = 698.25
$ws.Range("I68").Value = 698.25
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2094.75
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1283.75
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 698.25
$ws.Range("I71").Value = 698.25
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6284.25
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2228.25
$ws.Range("N71").ClearContents()

$ws.Range("H92").Value = 1296.5
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H120").Value = 32500
$ws.Range("I120").Value = 25000
$ws.Range("K120").Value = 75000
$ws.Range("M120").Value = -70162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 37748.5
$ws.Range("J57").Value = 37748.5
$ws.Range("L57").Value = 37748.5
$ws.Range("N57").Value = -39388.5

$ws.Range("H97").Value = 1218.5714
$ws.Range("I97").Value = 588.3333
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 588.3333
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -92.33330000000001
$ws.Range("N97").Value = -5992

$ws.Range("H102").Value = 863.6923
$ws.Range("I102").Value = 869.8
$ws.Range("J102").Value = 843.3333
$ws.Range("K102").Value = 869.8
$ws.Range("L102").Value = 843.3333
$ws.Range("M102").Value = 752.2
$ws.Range("N102").Value = -4087.3333

$ws.Range("H122").Value = 2527.3635
$ws.Range("I122").Value = 1734
$ws.Range("J122").Value = 3479.4
$ws.Range("K122").Value = 5202
$ws.Range("L122").Value = 10438.2
$ws.Range("M122").Value = -2752
$ws.Range("N122").Value = -15338.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H34").Value = 4000
$ws.Range("I34").Value = 4000
$ws.Range("K34").Value = 4000
$ws.Range("M34").Value = -3828

$ws.Range("H40").Value = 2507499.5
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 7499.75
$ws.Range("I46").Value = 4999.5
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 4999.5
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -4811.5
$ws.Range("N46").Value = -10376

$ws.Range("H82").Value = 1615
$ws.Range("I82").Value = 1486.6666
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1486.6666
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -1125.6666
$ws.Range("N82").Value = -2722

$ws.Range("H85").Value = 1615
$ws.Range("I85").Value = 1486.6666
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1486.6666
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -238.6666
$ws.Range("N85").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1707
$ws.Range("I122").Value = 1490.75
$ws.Range("K122").Value = 4472.25
$ws.Range("M122").Value = -2022.25

$ws.Range("H126").Value = 2253.6667
$ws.Range("I126").Value = 2034.75
$ws.Range("K126").Value = 6104.25
$ws.Range("M126").Value = -3634.25
